# "sua css man signup + navbar" - update In-charge / In-charge Full Name
# columns (E:F) for the signup-related rows (User Login / User Register /
# User Authorization) and the navbar-related rows (Home Page / Products
# List / Product Details / Cart Details) of the Product Backlog sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Backlog")

# Row 15 (User Authorization) had a stray formatting glitch: its E/F cells
# used the "odd" banding style instead of matching rows 13-14 (same
# HungND group). Copy the correct formats across before touching values.
$ws.Range("E13").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null
$ws.Range("F13").Copy() | Out-Null
$ws.Range("F15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Signup rows: User Login / User Register / User Authorization now also
# credit DucDT (back-end) alongside HungND (front-end). These cells carry a
# legacy "quote prefix" cell format (they store text that looks like it
# could be reinterpreted) - lead with an apostrophe so the text-literal
# formatting is preserved exactly as it was, matching the original file.
$ws.Range("E13").Value = "'HungND, DucDT"
$ws.Range("E14").Value = "'HungND, DucDT"
$ws.Range("E15").Value = "'HungND, DucDT"

# --- Navbar rows: Products List / Product Details / Cart Details now
# credit both AnhLH (front-end) and MinhVH (back-end).
$ws.Range("F10").Value = "Lê Hoàng Anh (FE-50), Vũ Nhật Minh (BE-50)"
$ws.Range("F11").Value = "Lê Hoàng Anh (FE-50), Vũ Nhật Minh (BE-50)"
$ws.Range("F12").Value = "Lê Hoàng Anh (FE-50), Vũ Nhật Minh (BE-50)"
$ws.Range("E11").Value = "AnhLH, MinhVH"
$ws.Range("E12").Value = "AnhLH, MinhVH"

$ws.Range("F13").Value = "Nguyễn Duy Hưng (FE-50), Đặng Trần Đức (BE-50)"
$ws.Range("F14").Value = "Nguyễn Duy Hưng (FE-50), Đặng Trần Đức (BE-50)"
$ws.Range("F15").Value = "Nguyễn Duy Hưng (FE-50), Đặng Trần Đức (BE-50)"

# --- Home Page now explicitly labelled as front-end work.
$ws.Range("F9").Value = "Lê Hoàng Anh (FE)"

# --- Restore the view: scrolled/selected around the signup rows.
$ws.Activate()
$ws.Range("F19").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
